$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = '''296.36'
$ws.Cells.Item(2,4).Style = 'Normal'
$ws.Cells.Item(2,5).Value = '''1.62%'
$ws.Cells.Item(2,5).Style = 'Normal'
$ws.Cells.Item(3,4).Value = '''41.93'
$ws.Cells.Item(3,4).Style = 'Normal'
$ws.Cells.Item(3,5).Value = '''3.03%'
$ws.Cells.Item(3,5).Style = 'Normal'
$ws.Cells.Item(4,4).Value = '''5.032'
$ws.Cells.Item(4,4).Style = 'Normal'
$ws.Cells.Item(4,5).Value = '''-0.11%'
$ws.Cells.Item(4,5).Style = 'Normal'
$ws.Cells.Item(5,4).Value = '''0.07576'
$ws.Cells.Item(5,4).Style = 'Normal'
$ws.Cells.Item(5,5).Value = '''2.80%'
$ws.Cells.Item(5,5).Style = 'Normal'
$ws.Cells.Item(6,2).Value = 'FTXToken'
$ws.Cells.Item(6,3).Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Cells.Item(6,4).Value = '''1.612'
$ws.Cells.Item(6,4).Style = 'Normal'
$ws.Cells.Item(6,5).Value = '''4.46%'
$ws.Cells.Item(6,5).Style = 'Normal'
$ws.Cells.Item(7,2).Value = 'MXToken'
$ws.Cells.Item(7,3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(7,4).Value = '''0.9312'
$ws.Cells.Item(7,4).Style = 'Normal'
$ws.Cells.Item(7,5).Value = '''0.50%'
$ws.Cells.Item(7,5).Style = 'Normal'
$ws.Cells.Item(8,2).Value = 'BTSEToken'
$ws.Cells.Item(8,3).Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Cells.Item(8,4).Value = '''2.408'
$ws.Cells.Item(8,4).Style = 'Normal'
$ws.Cells.Item(8,5).Value = '''3.44%'
$ws.Cells.Item(8,5).Style = 'Normal'
$ws.Cells.Item(9,2).Value = 'LiechtensteinCryptoassetsExchange'
$ws.Cells.Item(9,3).Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Cells.Item(9,4).Value = '''0.1198'
$ws.Cells.Item(9,4).Style = 'Normal'
$ws.Cells.Item(9,5).Value = '''3.94%'
$ws.Cells.Item(9,5).Style = 'Normal'
$ws.Cells.Item(10,2).Value = 'WazirX'
$ws.Cells.Item(10,3).Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Cells.Item(10,4).Value = '''0.1835'
$ws.Cells.Item(10,4).Style = 'Normal'
$ws.Cells.Item(10,5).Value = '''6.24%'
$ws.Cells.Item(10,5).Style = 'Normal'
$ws.Cells.Item(11,2).Value = 'MandalaExchangeToken'
$ws.Cells.Item(11,3).Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Cells.Item(11,4).Value = '''0.09005'
$ws.Cells.Item(11,4).Style = 'Normal'
$ws.Cells.Item(11,5).Value = '''4.80%'
$ws.Cells.Item(11,5).Style = 'Normal'
$ws.Cells.Item(12,2).Value = 'BitrueCoin'
$ws.Cells.Item(12,3).Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Cells.Item(12,4).Value = '''0.03974'
$ws.Cells.Item(12,4).Style = 'Normal'
$ws.Cells.Item(12,5).Value = '''-4.74%'
$ws.Cells.Item(12,5).Style = 'Normal'
$ws.Cells.Item(13,2).Value = 'BitMartToken'
$ws.Cells.Item(13,3).Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Cells.Item(13,4).Value = '''0.1051'
$ws.Cells.Item(13,4).Style = 'Normal'
$ws.Cells.Item(13,5).Value = '''-0.43%'
$ws.Cells.Item(13,5).Style = 'Normal'
$ws.Cells.Item(14,2).Value = 'BitForexToken'
$ws.Cells.Item(14,3).Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Cells.Item(14,4).Value = '''0.001280'
$ws.Cells.Item(14,4).Style = 'Normal'
$ws.Cells.Item(14,5).Value = '''1.08%'
$ws.Cells.Item(14,5).Style = 'Normal'
$ws.Cells.Item(15,2).Value = 'TigerCash'
$ws.Cells.Item(15,3).Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Cells.Item(15,4).Value = '''0.005829'
$ws.Cells.Item(15,4).Style = 'Normal'
$ws.Cells.Item(15,5).Value = '''-1.70%'
$ws.Cells.Item(15,5).Style = 'Normal'
$ws.Cells.Item(16,2).Value = 'LEO'
$ws.Cells.Item(16,3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(16,4).Value = '''3.365'
$ws.Cells.Item(16,4).Style = 'Normal'
$ws.Cells.Item(16,5).Value = '''-1.42%'
$ws.Cells.Item(16,5).Style = 'Normal'
$ws.Cells.Item(17,2).Value = 'GateToken'
$ws.Cells.Item(17,3).Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Cells.Item(17,4).Value = '''4.391'
$ws.Cells.Item(17,4).Style = 'Normal'
$ws.Cells.Item(17,5).Value = '''2.52%'
$ws.Cells.Item(17,5).Style = 'Normal'
$ws.Cells.Item(18,4).Value = '''0.3320'
$ws.Cells.Item(18,4).Style = 'Normal'
$ws.Cells.Item(18,5).Value = '''1.15%'
$ws.Cells.Item(18,5).Style = 'Normal'
$ws.Cells.Item(19,4).Value = '''7.892'
$ws.Cells.Item(19,4).Style = 'Normal'
$ws.Cells.Item(19,5).Value = '''2.20%'
$ws.Cells.Item(19,5).Style = 'Normal'
$ws.Cells.Item(20,4).Value = '''0.1419'
$ws.Cells.Item(20,4).Style = 'Normal'
$ws.Cells.Item(20,5).Value = '''2.93%'
$ws.Cells.Item(20,5).Style = 'Normal'
$ws.Cells.Item(21,4).Value = '''0.2998'
$ws.Cells.Item(21,4).Style = 'Normal'
$ws.Cells.Item(21,5).Value = '''4.04%'
$ws.Cells.Item(21,5).Style = 'Normal'
$ws.Cells.Item(22,4).Value = '''0.04065'
$ws.Cells.Item(22,4).Style = 'Normal'
$ws.Cells.Item(22,5).Value = '''5.48%'
$ws.Cells.Item(22,5).Style = 'Normal'
$ws.Cells.Item(23,5).Value = '''0.48%'
$ws.Cells.Item(23,5).Style = 'Normal'
$ws.Cells.Item(24,4).Value = '''0.003986'
$ws.Cells.Item(24,4).Style = 'Normal'
$ws.Cells.Item(24,5).Value = '''4.48%'
$ws.Cells.Item(24,5).Style = 'Normal'
$ws.Cells.Item(25,4).Value = '''0.0001230'
$ws.Cells.Item(25,4).Style = 'Normal'
$ws.Cells.Item(25,5).Value = '''-3.80%'
$ws.Cells.Item(25,5).Style = 'Normal'
$ws.Cells.Item(26,5).Value = '''0.11%'
$ws.Cells.Item(26,5).Style = 'Normal'
$ws.Cells.Item(38,4).Value = '''0.02417'
$ws.Cells.Item(38,4).Style = 'Normal'
$ws.Cells.Item(38,5).Value = '''3.71%'
$ws.Cells.Item(38,5).Style = 'Normal'
$ws.Cells.Item(39,4).Value = '''0.05206'
$ws.Cells.Item(39,4).Style = 'Normal'
$ws.Cells.Item(39,5).Value = '''4.30%'
$ws.Cells.Item(39,5).Style = 'Normal'
$ws.Cells.Item(40,4).Value = '''0.006058'
$ws.Cells.Item(40,4).Style = 'Normal'
$ws.Cells.Item(40,5).Value = '''6.58%'
$ws.Cells.Item(40,5).Style = 'Normal'
$ws.Cells.Item(41,4).Value = '''0.007781'
$ws.Cells.Item(41,4).Style = 'Normal'
$ws.Cells.Item(41,5).Value = '''1.40%'
$ws.Cells.Item(41,5).Style = 'Normal'
$ws.Cells.Item(42,4).Value = '''0.1334'
$ws.Cells.Item(42,4).Style = 'Normal'
$ws.Cells.Item(42,5).Value = '''4.02%'
$ws.Cells.Item(42,5).Style = 'Normal'
$ws.Cells.Item(43,4).Value = '''0.007534'
$ws.Cells.Item(43,4).Style = 'Normal'
$ws.Cells.Item(43,5).Value = '''2.56%'
$ws.Cells.Item(43,5).Style = 'Normal'
$ws.Cells.Item(44,4).Value = '''0.007225'
$ws.Cells.Item(44,4).Style = 'Normal'
$ws.Cells.Item(44,5).Value = '''1.83%'
$ws.Cells.Item(44,5).Style = 'Normal'
$ws.Cells.Item(45,4).Value = '''0.2980'
$ws.Cells.Item(45,4).Style = 'Normal'
$ws.Cells.Item(45,5).Value = '''-5.34%'
$ws.Cells.Item(45,5).Style = 'Normal'
$ws.Cells.Item(46,4).Value = '''0.00006773'
$ws.Cells.Item(46,4).Style = 'Normal'
$ws.Cells.Item(46,5).Value = '''5.57%'
$ws.Cells.Item(46,5).Style = 'Normal'
$ws.Cells.Item(47,4).Value = '''0.00000000750'
$ws.Cells.Item(47,4).Style = 'Normal'
$ws.Cells.Item(47,5).Value = '''0.02%'
$ws.Cells.Item(47,5).Style = 'Normal'
$ws.Cells.Item(48,4).Value = '''0.04278'
$ws.Cells.Item(48,4).Style = 'Normal'
$ws.Cells.Item(48,5).Value = '''148.94%'
$ws.Cells.Item(48,5).Style = 'Normal'
$ws.Cells.Item(49,4).Value = '''0.004199'
$ws.Cells.Item(49,4).Style = 'Normal'
$ws.Cells.Item(49,5).Value = '''-0.07%'
$ws.Cells.Item(49,5).Style = 'Normal'
$ws.Cells.Item(50,4).Value = '''0.00002099'
$ws.Cells.Item(50,4).Style = 'Normal'
$ws.Cells.Item(50,5).Value = '''0.02%'
$ws.Cells.Item(50,5).Style = 'Normal'
$ws.Cells.Item(51,4).Value = '''0.0001999'
$ws.Cells.Item(51,4).Style = 'Normal'
$ws.Cells.Item(51,5).Value = '''0.02%'
$ws.Cells.Item(51,5).Style = 'Normal'
